# Update "想去人数" (F column) values for the 漫展 rows that changed.
# The same updates apply to both the "展览" sheet and the "全部类型" sheet,
# which contain duplicated data.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 210
    3  = 253
    4  = 268
    5  = 803
    6  = 259
    7  = 6192
    8  = 45
    10 = 102
    11 = 66
    14 = 187
    15 = 445
    16 = 37
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
